# Edit script for "AMG Graphics Checklist for Submission 2021.docx"
# 1) Split the run "3 or more different shapes (line, rect, ellipse +++)" so that
#    "rect" is wrapped in a spellStart/spellEnd proofErr pair (as Word's proofer
#    would do when it flags "rect" as a possible misspelling).
# 2) Expand the contraction "it's" -> "it is" in the polyline/polygon description,
#    splitting that run into three runs.
#
# Because Range.InsertXML replaces the *entire* enclosing paragraph, each edit is
# performed by locating the full paragraph text with Find, then resupplying the
# complete paragraph (all of its original runs, with the target run split) as a
# Flat-OPC WordprocessingML fragment.

$d = $word.ActiveDocument

# --- Edit 1: "3 or more different shapes (line, rect, ellipse +++)" ---------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("“Live shapes” – 3 or more different shapes (line, rect, ellipse +++)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find the '3 or more different shapes' paragraph"
}
$xml1 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0915E345" w14:textId="35A88AFF" w:rsidR="004C7502" w:rsidRPr="00C86AB9" w:rsidRDefault="004C7502" w:rsidP="00CB7BB5"><w:r><w:t xml:space="preserve">“Live shapes” </w:t></w:r><w:r w:rsidR="00CE6776"><w:t>–</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00CE6776"><w:t xml:space="preserve">3 or more different shapes (line, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00CE6776"><w:t>rect</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00CE6776"><w:t>, ellipse +++)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng1.InsertXML($xml1)

# --- Edit 2: "...ensure it's at least within a 5-pixel radius..." -----------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("To close a polyline to form a polygon, repeat the steps above but for the final point ensure it’s at least within a 5-pixel radius of the starting point before releasing the mouse button.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the 'To close a polyline' paragraph"
}
$xml2 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5D45F005" w14:textId="4D930118" w:rsidR="004C7502" w:rsidRPr="00D1385D" w:rsidRDefault="007D536D" w:rsidP="00CB7BB5"><w:r><w:t xml:space="preserve">To close a polyline to form a polygon, repeat the steps above but for the final point ensure </w:t></w:r><w:r><w:t>it is</w:t></w:r><w:r><w:t xml:space="preserve"> at least within a 5-pixel radius of the starting point before releasing the mouse button.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng2.InsertXML($xml2)
